$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (row 11) Right count: 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (row 12) Right count: 69 -> 115
$ws.Range("B12").Value = 115

# Update the correct/total marks text: 65/84 -> 115/140
$ws.Range("E12").Value = "115/140"
